$d = $word.ActiveDocument

function Set-ParagraphText($index, $newText) {
    $para = $d.Paragraphs($index)
    $r = $para.Range
    $full = $d.Range($r.Start, $r.End - 1)
    # Force a real change so the engine rewrites the paragraph as a single
    # clean run (this also drops any stray w:proofErr spell-check markers).
    $full.Text = "`u{E000}PLACEHOLDER`u{E000}"
    $para2 = $d.Paragraphs($index)
    $r2 = $para2.Range
    $full2 = $d.Range($r2.Start, $r2.End - 1)
    $full2.Text = $newText
}

# --- 1) " do teclado, como se fosse renomear uma pasta no windows." ---
Set-ParagraphText 9 "Podemos reescrever em uma célula clicando duas vezes com o botão esquerdo ou apertando o botão f2 do teclado, como se fosse renomear uma pasta no windows."

# --- 2) "... igual estavamos fazendo antes ..." ---
Set-ParagraphText 17 "Quando criamos essa tabela, não podemos mais copiar linhas inteiras de células e colar abaixo igual estavamos fazendo antes, mas, para mudar isso, podemos clicar dentro da nossa tabela, ir para a guia de Design e clicar em “converter em intervalo”. Após clicar em sim a sua tabela continua formatada apesar de ter se tornado intervalo novamente."

# --- 3) "... ja que o excel pega do windows a informacao ..." ---
Set-ParagraphText 31 "Quando colocamos valores de moedas, não precisamos colocar cifrão ,0 nem nada. Basta digitar os valores, selecionar todas as células e na guia da página inicial e no bloco número, temos uma opção de “Formato de Número de Contabilização”. Ao clicar no botão, todas as células serão automaticamente configuradas para reais Brasil, já que o excel pega do windows a informação de que estamos no Brasil."

# --- 4) "... guia da Pagina inicial ..." ---
Set-ParagraphText 33 "Para fazer uma soma total de uma coluna, basta clicar na célula onde ficará o resultado e depois no ícone de soma no último bloco da guia da Pagina inicial, ou apertar ALT+=, com ela selecionada. Em seguida, escolher quais são as células que serão somadas. Pode ser feita por seleção manual ou apenas digitar de qual a qual é (C3:C19 [indica que é pra somar de c3 até c19])."

# --- 5) "... o proprio excel faz." ---
Set-ParagraphText 37 "Fórmula é quando escrevemos na mão a célula e operadores e a função é algo que o próprio excel faz."

# --- 6) "... que e excel faz a conta ..." ---
Set-ParagraphText 38 "A vantagem de utilizar fórmulas e funções é que não temos a necessidade de ficar alterando os valores resultados na célula sempre que mudarmos algo de um estoque, por exemplo. Basta mudar a quantidade e/ou o preço que e excel faz a conta e nos mostra o resultado final novo imediatamente."

# --- 7) Append new content at the very end of the document (Aula 4 section) ---

# Paragraph 42 currently holds a single space; extend it with the new
# sentence about selecting non-adjacent ranges (same run formatting, so it
# naturally collapses into a single run, matching the source paragraph's
# leading " " + new sentence).
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastIndex)
$lastRng = $lastPara.Range
$lastFull = $d.Range($lastRng.Start, $lastRng.End - 1)
$lastFull.Text = " Podemos selecionar intervalo de dados não adjacentes. Primeiro clica e seleciona o primeiro intervalo que deseja, depois aperta o CTRL e clica e arrasta selecionando o segundo intervalo que dejesa."

function Add-ListParagraph($text, $ilvl) {
    # ListLevelNumber is 1-based (ilvl 0 => level 1, ilvl 1 => level 2, ...)
    $level = $ilvl + 1
    $count = $d.Paragraphs.Count
    $prevPara = $d.Paragraphs($count)
    $prevRng = $prevPara.Range
    $prevRng.InsertParagraphAfter()
    $newCount = $d.Paragraphs.Count
    $newPara = $d.Paragraphs($newCount)
    if ($newPara.Range.ListFormat.ListLevelNumber -ne $level) {
        $newPara.Range.ListFormat.ListLevelNumber = $level
    }
    $newRng = $newPara.Range
    $newFull = $d.Range($newRng.Start, $newRng.End - 1)
    $newFull.Text = $text
}

Add-ListParagraph "Para criar um gráfico basta selecionar os dois intervalos de interesse e depois ir na guia “inserir”, bloco “gráficos” e clicar no botão “gráficos de colunas ou de barras” ." 1
Add-ListParagraph "Ao clicar no gráfico surge uma nova guia que não estava lá antes, para podermos mexer no “design de gráfico”, igual ao que aconteceu com a tabela." 1
Add-ListParagraph "Nessa guia podemos clicar em “Mover Gráfico”, selecionar a opção nova planilha, na janela que apareceu, e digitar um novo nome." 2
Add-ListParagraph "Essa estratégia cria uma nova planilha para o gráfico que estamos clicados." 2
Add-ListParagraph "Nós podemos fazer alterações no gráfico na guia da página inicial igual ao que fizemos na planilha, mas, caso mesmo assim o gráfico não tenha ficado bom, podemos “Alterar tipo de gráfico” na guia de Design do gráfico, escolhendo assim algum que fique mais proporcional e melhor." 1
Add-ListParagraph "Os dois tipos de gráficos mais utilizados são o “Barras e colunas” e o “gráfico de pizza”." 1
Add-ListParagraph "Gráfico de pizza coloca os percentuais automaticamente." 1

Write-Output "done"
